$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "See full schedule"
$ws.Range("B2").Value = 45447
$ws.Range("C2").Value = "Beyond the headlines to the heart of the news of the day. Al Jazeera gets the Inside Story from some of the best minds in the Middle East and beyond."
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = $false

$ws.Range("A3").Value = "Migration | Today's latest from Al Jazeera"
$ws.Range("B3").Value = 45447
$ws.Range("C3").Value = "Order will allow US to shut off asylum requests and deny entry to migrants once daily threshold met, US media reporting. Published On 3 Jun 20243 Jun"
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = $false

$ws.Range("A4").Value = "Latin America News | Today's latest from Al Jazeera"
$ws.Range("B4").Value = 45445
$ws.Range("C4").Value = "Sheinbaum's resounding win offers a bright spot for the Latin American left amid a string of setbacks.   Claudia Sheinbaum pumps her fist as she acknowledges"
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = $false

$ws.Range("A5").Value = "Leaders praise 'historic' victory as Sheinbaum triumphs in Mexican ..."
$ws.Range("B5").Value = 45447
$ws.Range("C5").Value = "By contrast, in countries like Colombia and Chile, left-leaning presidents have seen their popularity ebb as they struggle to make progress on"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = $false

$ws.Range("A6").Value = "Mexico's election puts Lopez Obrador's stance on Israel under ..."
$ws.Range("B6").Value = 45445
$ws.Range("C6").Value = "In Colombia, meanwhile, Gustavo Petro became the first left-wing leader to win the modern presidency."
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = $false

$ws.Range("A7").Value = "Brazil withdraws ambassador to Israel after Gaza war criticism ..."
$ws.Range("B7").Value = 45441
$ws.Range("C7").Value = "Colombian President Gustavo Petro, who has also severed ties with Israel. Both Brazil and Colombia have supported South Africa's complaint"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = $false

$ws.Range("A8").Value = "Colombia lawmakers pass bullfighting ban | News | Al Jazeera"
$ws.Range("B8").Value = 45441
$ws.Range("C8").Value = "Colombia's Congress has passed legislation banning bullfighting. Lawmakers passed the bill 93-2 on Tuesday. Activists have spent many years"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = $false

$ws.Range("A9").Value = "Colombia's ex-President Uribe charged with witness tampering ..."
$ws.Range("B9").Value = 45436
$ws.Range("C9").Value = "Former Colombia President Alvaro Uribe. Uribe, who was president from 2002 to 2010, has denied any wrongdoing and has accused Colombia's"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = $false

$ws.Range("A10").Value = "'Major non-NATO ally': What does Biden's new Kenya pledge mean ..."
$ws.Range("B10").Value = 45435
$ws.Range("C10").Value = "The US has currently designated 18 countries as MNNAs. These include Argentina, Australia, Bahrain, Brazil, Colombia, Egypt, Israel, Japan,"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = $false

$ws.Range("A11").Value = "Why is Israel angry some EU countries are recognising Palestine ..."
$ws.Range("B11").Value = 45434
$ws.Range("C11").Value = "Colombia's Petro orders opening of embassy in West Bank's Ramallah. Petro had recalled the Colombian ambassador from Tel Aviv, where the"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = $false

$ws.Range("A12").Value = "Colombia's Petro orders opening of embassy in West Bank's ..."
$ws.Range("B12").Value = 45434
$ws.Range("C12").Value = "Colombian President Gustavo Petro has ordered the opening of an embassy in the Palestinian city of Ramallah, Foreign Minister Luis Gilberto"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = $false

$ws.Range("A13").Value = "The ICC is not in the business of peacemaking, but it can deliver ..."
$ws.Range("B13").Value = 45434
$ws.Range("C13").Value = "But more `"complicated`" peace negotiations do not necessarily mean `"worse`" peace negotiations. Take Colombia, for example, where the ICC had a"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = $false

$ws.Range("A14").Value = "Mapping which countries recognise Palestine in 2024 | Israel ..."
$ws.Range("B14").Value = 45433
$ws.Range("C14").Value = "Colombia. 2015: Saint Lucia. 2014: Sweden. 2013: Guatemala, Haiti, the Vatican. 2012: Thailand. 2011: Chile, Guyana, Peru, Suriname, Uruguay"
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = $false

$ws.Range("A15").Value = "Are seed-sowing drones the answer to global deforestation ..."
$ws.Range("B15").Value = 45429
$ws.Range("C15").Value = "And, in Colombia, internal violence and displacement have pushed armed groups, farmers and cattle farmers into the forests, causing more"
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = $false

$ws.Range("A16").Value = "Colombia hunts for assailants after Bogota prison director shot dead ..."
$ws.Range("B16").Value = 45428
$ws.Range("C16").Value = "The new director of one of Colombia's biggest prisons has been shot dead, the authorities said, after receiving threats against him and his"
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = $false

$ws.Range("A17").Value = "Children of the Darien Gap | Migration | Al Jazeera"
$ws.Range("B17").Value = 45419
$ws.Range("C17").Value = "Colombia and Panama that is the only land route for migrants heading north from South America. Amid historic regional migration and new"
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = $false

$ws.Range("A18").Value = "Colombia president cuts ties with Israel over war on Gaza | Israel ..."
$ws.Range("B18").Value = 45413
$ws.Range("C18").Value = "Colombia's President Gustavo Petro says the country is cutting diplomatic ties with Israel over its war on Gaza."
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = $false

$ws.Range("A19").Value = "Arrests at Columbia University as New York City police clear Gaza ..."
$ws.Range("B19").Value = 45412
$ws.Range("C19").Value = "Violent clashes erupt at UCLA between pro-Palestinian protesters and pro-Israeli counter-demonstrators."
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = $false

$ws.Range("A20").Value = "Colombia to cut diplomatic ties with Israel over Gaza war, Petro says ..."
$ws.Range("B20").Value = 45412
$ws.Range("C20").Value = "Colombia to cut diplomatic ties with Israel over Gaza war, Petro says. Colombian President Gustavo Petro, a staunch critic of Israel's war in"
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = $false

$ws.Range("A21").Value = "Huge crowds protest Colombian president's planned reforms ..."
$ws.Range("B21").Value = 45403
$ws.Range("C21").Value = "Huge crowds protest Colombian president's planned reforms. Protesters call Gustavo Petro's policies 'dire' as his government attempts to reform"
$ws.Range("E21").Value = 2
$ws.Range("F21").Value = $false

$ws.Range("A22").Value = "Former Colombian President Alvaro Uribe blasts impending criminal ..."
$ws.Range("B22").Value = 45391
$ws.Range("C22").Value = "Prosecutors announced this week they intended to pursue the right-wing Uribe on charges of witness tampering and fraud."
$ws.Range("E22").Value = 1
$ws.Range("F22").Value = $false

$ws.Range("A23").Value = "Forced from home, these Colombians struggle to live in a basketball ..."
$ws.Range("B23").Value = 45391
$ws.Range("C23").Value = "Forced from home, Colombians build a life in a basketball stadium · An invisible crisis · A sanctuary from gunfire · Resources stretched thin."
$ws.Range("E23").Value = 2
$ws.Range("F23").Value = $false

$ws.Range("A24").Value = "Colombia seeks to join Gaza genocide case against Israel at ICJ ..."
$ws.Range("B24").Value = 45386
$ws.Range("C24").Value = "Colombia seeks to join Gaza genocide case against Israel at ICJ. Bogota calls on the World Court to ensure 'the safety' and 'the very existence"
$ws.Range("E24").Value = 2
$ws.Range("F24").Value = $false

$ws.Range("A25").Value = "Colombia and Panama failing to protect migrants in Darien Gap ..."
$ws.Range("B25").Value = 45384
$ws.Range("C25").Value = "In a report on Wednesday, the rights group said the Colombian and Panamanian authorities have not protected people transiting through the Darien"
$ws.Range("E25").Value = 2
$ws.Range("F25").Value = $false

$ws.Range("A26").Value = "Colombia expels Argentina's diplomats after Milei calls Petro ..."
$ws.Range("B26").Value = 45378
$ws.Range("C26").Value = "Colombia expels Argentina's diplomats after Milei calls Petro 'terrorist'. Argentina's president calls his Colombian counterpart a 'terrorist',"
$ws.Range("E26").Value = 3
$ws.Range("F26").Value = $false

$ws.Range("A27").Value = "Tonnes of cocaine seized after high-speed boat chase in Colombia ..."
$ws.Range("B27").Value = 45377
$ws.Range("C27").Value = "Over five tonnes of cocaine have been seized in Colombia after two drug bust operations that involved a boat chase."
$ws.Range("E27").Value = 2
$ws.Range("F27").Value = $false

$ws.Range("A28").Value = "In Colombia, hunting poachers, not drug traffickers | Wildlife News ..."
$ws.Range("B28").Value = 45371
$ws.Range("C28").Value = "One of the most biodiverse countries in the world, Colombia is increasingly vulnerable to illicit wildlife trafficking; in 2023, the wildlife"
$ws.Range("E28").Value = 2
$ws.Range("F28").Value = $false

$ws.Range("A29").Value = "Colombia names attorney general amid political unrest under ..."
$ws.Range("B29").Value = 45362
$ws.Range("C29").Value = "Colombia names attorney general amid political unrest under Gustavo Petro. Luz Adriana Camargo Garzon will lead probes into President Petro and"
$ws.Range("E29").Value = 2
$ws.Range("F29").Value = $false

$ws.Range("A30:F35").EntireRow.Delete()

Write-Host "Done updating rows and deleting 30-35"